$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 20833530
$ws.Range("I8").Value = 47619124
$ws.Range("J8").Value = 287.8889
$ws.Range("K8").Value = 142857372
$ws.Range("L8").Value = 863.6667
$ws.Range("M8").Value = -142857233
$ws.Range("N8").Value = -1141.6667

# Row 9
$ws.Range("H9").Value = 230.75
$ws.Range("J9").Value = 6.25
$ws.Range("L9").Value = 6.25
$ws.Range("N9").Value = -344.25

# Row 11
$ws.Range("H11").Value = 80.23077000000001
$ws.Range("I11").Value = 80.23077000000001
$ws.Range("K11").Value = 80.23077000000001
$ws.Range("M11").Value = 59.76922999999999

# Row 28
$ws.Range("H28").Value = 267.86667
$ws.Range("I28").Value = 306.63635
$ws.Range("K28").Value = 306.63635
$ws.Range("M28").Value = 178.36365

# Row 64
$ws.Range("H64").Value = 103332.5
$ws.Range("I64").Value = 200000
$ws.Range("K64").Value = 200000
$ws.Range("M64").Value = -199752

# Row 67
$ws.Range("H67").Value = 103332.5
$ws.Range("I67").Value = 200000
$ws.Range("K67").Value = 200000
$ws.Range("M67").Value = -199142

# Row 76
$ws.Range("H76").Value = 4560.9644
$ws.Range("I76").Value = 3846
$ws.Range("J76").Value = 5385.923
$ws.Range("K76").Value = 3846
$ws.Range("L76").Value = 5385.923
$ws.Range("M76").Value = -3531
$ws.Range("N76").Value = -6015.923

# Row 79
$ws.Range("H79").Value = 4560.9644
$ws.Range("I79").Value = 3846
$ws.Range("J79").Value = 5385.923
$ws.Range("K79").Value = 3846
$ws.Range("L79").Value = 5385.923
$ws.Range("M79").Value = -2754
$ws.Range("N79").Value = -7569.923

# Row 132
$ws.Range("H132").Value = 1726888.2
$ws.Range("I132").Value = 2799.6345
$ws.Range("K132").Value = 8398.9035
$ws.Range("M132").Value = -5868.9035

# Row 137
$ws.Range("H137").Value = 7735.4224
$ws.Range("I137").Value = 10134.322
$ws.Range("K137").Value = 30402.966
$ws.Range("M137").Value = -27852.966

# Row 138
$ws.Range("H138").Value = 283287
$ws.Range("J138").Value = 5089.25
$ws.Range("L138").Value = 15267.75
$ws.Range("N138").Value = -25547.75


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5668.476
$ws.Range("I32").Value = 5495.807
$ws.Range("K32").Value = 5495.807
$ws.Range("M32").Value = -5208.807

# Row 45
$ws.Range("H45").Value = 6657.077
$ws.Range("J45").Value = 7402.8
$ws.Range("L45").Value = 7402.8
$ws.Range("N45").Value = -8156.8

# Row 74
$ws.Range("H74").Value = 1646.7931
$ws.Range("I74").Value = 932.7646999999999
$ws.Range("J74").Value = 2658.3333
$ws.Range("K74").Value = 932.7646999999999
$ws.Range("L74").Value = 2658.3333
$ws.Range("M74").Value = -58.76469999999995
$ws.Range("N74").Value = -4406.3333

# Row 77
$ws.Range("H77").Value = 1646.7931
$ws.Range("I77").Value = 932.7646999999999
$ws.Range("J77").Value = 2658.3333
$ws.Range("K77").Value = 4663.8235
$ws.Range("L77").Value = 13291.6665
$ws.Range("M77").Value = -295.8234999999995
$ws.Range("N77").Value = -22027.6665

# Row 132
$ws.Range("H132").Value = 4006.0715
$ws.Range("J132").Value = 6456.857
$ws.Range("L132").Value = 19370.571
$ws.Range("N132").Value = -24430.571


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2707.8276
$ws.Range("I20").Value = 2039.9412
$ws.Range("J20").Value = 3654
$ws.Range("K20").Value = 2039.9412
$ws.Range("L20").Value = 3654
$ws.Range("M20").Value = -1792.9412
$ws.Range("N20").Value = -4148

# Row 69
$ws.Range("L69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("L72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("N72").ClearContents()

# Row 107
$ws.Range("H107").Value = 2024.8667
$ws.Range("I107").Value = 2014.4166
$ws.Range("J107").Value = 2066.6667
$ws.Range("K107").Value = 2014.4166
$ws.Range("L107").Value = 2066.6667
$ws.Range("M107").Value = -94.41660000000002
$ws.Range("N107").Value = -5906.6667

# Row 132
$ws.Range("H132").Value = 76999
$ws.Range("J132").Value = 76999
$ws.Range("L132").Value = 76999
$ws.Range("N132").Value = -87119

# Row 134
$ws.Range("H134").Value = 5260.278
$ws.Range("I134").Value = 5273.1113
$ws.Range("J134").Value = 5221.778
$ws.Range("K134").Value = 15819.3339
$ws.Range("L134").Value = 15665.334
$ws.Range("M134").Value = -13284.3339
$ws.Range("N134").Value = -20735.334

# Row 135
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 525.8333
$ws.Range("I22").Value = 525.8333
$ws.Range("K22").Value = 525.8333
$ws.Range("M22").Value = -175.8333

# Row 31
$ws.Range("H31").Value = 8755.299999999999
$ws.Range("I31").Value = 10985.125
$ws.Range("J31").Value = 6206.9287
$ws.Range("K31").Value = 10985.125
$ws.Range("L31").Value = 6206.9287
$ws.Range("M31").Value = -10690.125
$ws.Range("N31").Value = -6796.9287

# Row 34
$ws.Range("H34").Value = 8755.299999999999
$ws.Range("I34").Value = 10985.125
$ws.Range("J34").Value = 6206.9287
$ws.Range("K34").Value = 10985.125
$ws.Range("L34").Value = 6206.9287
$ws.Range("M34").Value = -10783.125
$ws.Range("N34").Value = -6610.9287

# Row 86
$ws.Range("H86").Value = 5763.1577
$ws.Range("I86").Value = 5418.6924
$ws.Range("J86").Value = 6509.5
$ws.Range("K86").Value = 5418.6924
$ws.Range("L86").Value = 6509.5
$ws.Range("M86").Value = -4295.6924
$ws.Range("N86").Value = -8755.5

# Row 89
$ws.Range("H89").Value = 5763.1577
$ws.Range("I89").Value = 5418.6924
$ws.Range("J89").Value = 6509.5
$ws.Range("K89").Value = 27093.462
$ws.Range("L89").Value = 32547.5
$ws.Range("M89").Value = -21477.462
$ws.Range("N89").Value = -43779.5

# Row 132
$ws.Range("H132").Value = 1336.091
$ws.Range("I132").Value = 1207
$ws.Range("K132").Value = 3621
$ws.Range("M132").Value = -1091

# Row 134
$ws.Range("H134").Value = 3863.8
$ws.Range("I134").Value = 2836.8823
$ws.Range("K134").Value = 8510.6469
$ws.Range("M134").Value = -5975.6469

# Row 141
$ws.Range("H141").Value = 184931.36
$ws.Range("I141").Value = 70000
$ws.Range("J141").Value = 197245.42
$ws.Range("K141").Value = 70000
$ws.Range("L141").Value = 197245.42
$ws.Range("M141").Value = -64820
$ws.Range("N141").Value = -207605.42


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 1997
$ws.Range("J25").Value = 2330.5881
$ws.Range("L25").Value = 6991.7643
$ws.Range("N25").Value = -7329.7643

# Row 30
$ws.Range("H30").Value = 1997
$ws.Range("J30").Value = 2330.5881
$ws.Range("L30").Value = 6991.7643
$ws.Range("N30").Value = -7195.7643

# Row 34
$ws.Range("H34").Value = 992.125
$ws.Range("I34").Value = 862.5
$ws.Range("K34").Value = 2587.5
$ws.Range("M34").Value = -2503.5

# Row 39
$ws.Range("H39").Value = 482.8889
$ws.Range("I39").Value = 287.7647
$ws.Range("J39").Value = 3800
$ws.Range("K39").Value = 863.2941000000001
$ws.Range("L39").Value = 11400
$ws.Range("M39").Value = -569.2941000000001
$ws.Range("N39").Value = -11988

# Row 54
$ws.Range("H54").Value = 928.5714
$ws.Range("I54").Value = 583.3333
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 1749.9999
$ws.Range("L54").Value = 9000
$ws.Range("M54").Value = -1190.9999
$ws.Range("N54").Value = -10118

# Row 55
$ws.Range("H55").Value = 10300.556
$ws.Range("J55").Value = 12772.714
$ws.Range("L55").Value = 38318.142
$ws.Range("N55").Value = -38672.142

# Row 97
$ws.Range("M97").Value = -183347.6
$ws.Range("H97").Value = 61281.2
$ws.Range("I97").Value = 61281.2
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 183843.6
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 98
$ws.Range("H98").Value = 2531.125
$ws.Range("J98").Value = 2375
$ws.Range("L98").Value = 7125
$ws.Range("N98").Value = -10121


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 24141.143
$ws.Range("I80").Value = 27332
$ws.Range("J80").Value = 4996
$ws.Range("K80").Value = 27332
$ws.Range("L80").Value = 4996
$ws.Range("M80").Value = -26334
$ws.Range("N80").Value = -6992

# Row 83
$ws.Range("H83").Value = 24141.143
$ws.Range("I83").Value = 27332
$ws.Range("J83").Value = 4996
$ws.Range("K83").Value = 136660
$ws.Range("L83").Value = 24980
$ws.Range("M83").Value = -131668
$ws.Range("N83").Value = -34964

# Row 122
$ws.Range("H122").Value = 10526.143
$ws.Range("I122").Value = 7759.048
$ws.Range("J122").Value = 18827.428
$ws.Range("K122").Value = 23277.144
$ws.Range("L122").Value = 56482.284
$ws.Range("M122").Value = -20827.144
$ws.Range("N122").Value = -61382.284

# Row 132
$ws.Range("H132").Value = 1842.3864
$ws.Range("I132").Value = 1835.619
$ws.Range("J132").Value = 1984.5
$ws.Range("K132").Value = 5506.857
$ws.Range("L132").Value = 5953.5
$ws.Range("M132").Value = -2976.857
$ws.Range("N132").Value = -11013.5


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1060.1025
$ws.Range("I16").Value = 1094.2903
$ws.Range("J16").Value = 927.625
$ws.Range("K16").Value = 1094.2903
$ws.Range("L16").Value = 927.625
$ws.Range("M16").Value = -924.2902999999999
$ws.Range("N16").Value = -1267.625

# Row 82
$ws.Range("H82").Value = 1865.1428
$ws.Range("I82").Value = 1514.5
$ws.Range("J82").Value = 2332.6667
$ws.Range("K82").Value = 1514.5
$ws.Range("L82").Value = 2332.6667
$ws.Range("M82").Value = -1153.5
$ws.Range("N82").Value = -3054.6667

# Row 85
$ws.Range("H85").Value = 1865.1428
$ws.Range("I85").Value = 1514.5
$ws.Range("J85").Value = 2332.6667
$ws.Range("K85").Value = 1514.5
$ws.Range("L85").Value = 2332.6667
$ws.Range("M85").Value = -266.5
$ws.Range("N85").Value = -4828.6667

# Row 132
$ws.Range("H132").Value = 499248.78
$ws.Range("I132").Value = 994566.3
$ws.Range("K132").Value = 2983698.9
$ws.Range("M132").Value = -2981168.9


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 22037.305
$ws.Range("I126").Value = 32007.467
$ws.Range("K126").Value = 96022.401
$ws.Range("M126").Value = -93552.401

